$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.837.41'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '1.936.11'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''243.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '''0.4911'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').Value = '''0.2960'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '''0.06901'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').Value = '''19.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').Value = '''104.89'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.02%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.943.75'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '''0.07793'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('D14').Value = '''5.355'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.99%  '
$ws.Range('D15').Value = '''0.7016'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '''273.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.99%  '
$ws.Range('D17').Value = '30.843.66'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').Value = '''0.000007736'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '''13.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''5.602'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '''1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '2.195.50'
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').Value = '''1.003'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '''6.551'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('D25').Value = '''9.867'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').Value = '''165.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.15%  '
$ws.Range('D27').Value = '''19.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').Value = '''2.160'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('E29').Value = '  -1.32%  '
$ws.Range('D30').Value = '''1.392'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('D31').Value = '''1.559'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('D32').Value = '''4.573'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('D33').Value = '''4.383'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.58%  '
$ws.Range('D34').Value = '''0.04900'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').Value = '''0.7628'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = '''2.738'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('D39').Value = '''0.02013'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').Value = '''79.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.34%  '
$ws.Range('D41').Value = '''2.662'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range('D42').Value = '''6.526'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('D43').Value = '''2.085'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range('D44').Value = '''0.9079'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').Value = '''0.4452'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('D46').Value = '''107.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''7.762'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.89%  '
$ws.Range('D49').Value = '''995.84'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('D51').Value = '''36.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.37%  '
